$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 374
$ws.Range("I2").Value = 1160
$ws.Range("J2").Value = 4556
$ws.Range("K2").Value = 19
$ws.Range("L2").Value = 1288
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 837
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 10
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 55
$ws.Range("S2").Value = 522
$ws.Range("T2").Value = 780
$ws.Range("U2").Value = 62
$ws.Range("V2").Value = 7157
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 7055
$ws.Range("Y2").Value = 13
$ws.Range("Z2").Value = 103
$ws.Range("AA2").Value = 47
